$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 351818
$ws.Range("D2").Value = 446868126
$ws.Range("C8").Value = 927
$ws.Range("D8").Value = 1363596
$ws.Range("C10").Value = 125768
$ws.Range("D10").Value = 184173454
$ws.Range("C12").Value = 65999
$ws.Range("D12").Value = 95262366
$ws.Range("C16").Value = 4161
$ws.Range("D16").Value = 5909075
$ws.Range("C21").Value = 8709
$ws.Range("D21").Value = 12152134
$ws.Range("C23").Value = 83750
$ws.Range("D23").Value = 103932946
$ws.Range("C29").Value = 34420
$ws.Range("D29").Value = 50369864
$ws.Range("C32").Value = 12523
$ws.Range("D32").Value = 18016746
$ws.Range("C37").Value = 2301
$ws.Range("D37").Value = 3241261
$ws.Range("C38").Value = 104656
$ws.Range("D38").Value = 131033102
$ws.Range("C46").Value = 46916
$ws.Range("D46").Value = 68717144
$ws.Range("C48").Value = 9986
$ws.Range("D48").Value = 14312970
$ws.Range("C53").Value = 2909
$ws.Range("D53").Value = 4078593
$ws.Range("C54").Value = 74761
$ws.Range("D54").Value = 93487820
$ws.Range("C61").Value = 30176
$ws.Range("D61").Value = 44214496
$ws.Range("C64").Value = 12207
$ws.Range("D64").Value = 17635115
$ws.Range("C70").Value = 1869
$ws.Range("D70").Value = 2625827
$ws.Range("C72").Value = 22513
$ws.Range("D72").Value = 29401420
$ws.Range("C73").Value = 37
$ws.Range("D73").Value = 48419
$ws.Range("C76").Value = 8293
$ws.Range("D76").Value = 12143124
$ws.Range("C78").Value = 5685
$ws.Range("D78").Value = 8266931
$ws.Range("C79").Value = 560
$ws.Range("D79").Value = 790593
$ws.Range("C80").Value = 362
$ws.Range("D80").Value = 510948
$ws.Range("C81").Value = 152908
$ws.Range("D81").Value = 189785875
$ws.Range("C87").Value = 67596
$ws.Range("D87").Value = 99013356
$ws.Range("C90").Value = 32180
$ws.Range("D90").Value = 46554234
$ws.Range("C92").Value = 2991
$ws.Range("D92").Value = 4315182
$ws.Range("C94").Value = 3543
$ws.Range("D94").Value = 5000101
$ws.Range("C95").Value = 37861
$ws.Range("D95").Value = 51354821
$ws.Range("C99").Value = 9396
$ws.Range("D99").Value = 13804262
$ws.Range("C101").Value = 8734
$ws.Range("D101").Value = 12673426
$ws.Range("C103").Value = 599
$ws.Range("D103").Value = 846930
$ws.Range("C105").Value = 16556
$ws.Range("D105").Value = 30986147
$ws.Range("C108").Value = 3839
$ws.Range("D108").Value = 7702389
$ws.Range("C110").Value = 5451
$ws.Range("D110").Value = 11055384
$ws.Range("C112").Value = 254
$ws.Range("D112").Value = 510930
$ws.Range("C115").Value = 153885
$ws.Range("D115").Value = 189917997
$ws.Range("C121").Value = 56584
$ws.Range("D121").Value = 82859217
$ws.Range("C122").Value = 102
$ws.Range("D122").Value = 151459
$ws.Range("C123").Value = 30443
$ws.Range("D123").Value = 44098626
$ws.Range("C127").Value = 2838
$ws.Range("D127").Value = 3998497
$ws.Range("C129").Value = 615399
$ws.Range("D129").Value = 812273807
$ws.Range("C134").Value = 1596
$ws.Range("D134").Value = 2363603
$ws.Range("C136").Value = 237056
$ws.Range("D136").Value = 348238519
$ws.Range("C137").Value = 540
$ws.Range("D137").Value = 804647
$ws.Range("C139").Value = 221230
$ws.Range("D139").Value = 321695694
$ws.Range("C142").Value = 3021
$ws.Range("D142").Value = 4246734
$ws.Range("C145").Value = 8515
$ws.Range("D145").Value = 11996993
$ws.Range("C148").Value = 49566
$ws.Range("D148").Value = 66084078
$ws.Range("C154").Value = 15283
$ws.Range("D154").Value = 22398082
$ws.Range("C155").Value = 4151
$ws.Range("D155").Value = 5990826
$ws.Range("C161").Value = 19818
$ws.Range("D161").Value = 26199666
$ws.Range("C165").Value = 8204
$ws.Range("D165").Value = 11935032
$ws.Range("C167").Value = 5824
$ws.Range("D167").Value = 8383492
$ws.Range("C172").Value = 28771
$ws.Range("D172").Value = 57959648
$ws.Range("C178").Value = 94645
$ws.Range("D178").Value = 117830425
$ws.Range("C185").Value = 36174
$ws.Range("D185").Value = 53015854
$ws.Range("C187").Value = 14394
$ws.Range("D187").Value = 20792331
$ws.Range("C191").Value = 2046
$ws.Range("D191").Value = 2875484
$ws.Range("C193").Value = 258123
$ws.Range("D193").Value = 319674038
$ws.Range("C195").Value = 187
$ws.Range("D195").Value = 267891
$ws.Range("C201").Value = 92394
$ws.Range("D201").Value = 135362452
$ws.Range("C204").Value = 36289
$ws.Range("D204").Value = 52245327
$ws.Range("C207").Value = 5477
$ws.Range("D207").Value = 7799272
$ws.Range("C210").Value = 6076
$ws.Range("D210").Value = 8412405
$ws.Range("C213").Value = 287156
$ws.Range("D213").Value = 354249549
$ws.Range("C220").Value = 657
$ws.Range("D220").Value = 956370
$ws.Range("C222").Value = 102196
$ws.Range("D222").Value = 149446849
$ws.Range("C225").Value = 56704
$ws.Range("D225").Value = 81912028
$ws.Range("C228").Value = 4938
$ws.Range("D228").Value = 6926402
$ws.Range("C231").Value = 7541
$ws.Range("D231").Value = 10437508
$ws.Range("C234").Value = 115784
$ws.Range("D234").Value = 144216629
$ws.Range("C241").Value = 52697
$ws.Range("D241").Value = 77170558
$ws.Range("C243").Value = 13905
$ws.Range("D243").Value = 20003147
$ws.Range("C245").Value = 1965
$ws.Range("D245").Value = 2815482
$ws.Range("C247").Value = 3096
$ws.Range("D247").Value = 4331017
$ws.Range("C248").Value = 283783
$ws.Range("D248").Value = 357735766
$ws.Range("C249").Value = 197
$ws.Range("D249").Value = 245828
$ws.Range("C257").Value = 104566
$ws.Range("D257").Value = 153167368
$ws.Range("C258").Value = 91
$ws.Range("D258").Value = 131564
$ws.Range("C259").Value = 7
$ws.Range("D259").Value = 10500
$ws.Range("C260").Value = 73469
$ws.Range("D260").Value = 106476405
$ws.Range("C262").Value = 2557
$ws.Range("D262").Value = 3599706
$ws.Range("C264").Value = 9
$ws.Range("D264").Value = 13500
$ws.Range("C265").Value = 6015
$ws.Range("D265").Value = 8434539
